# Commit: "stateless entities outside the US"
#
# The region table has three repeated 8-column blocks (B:I = M_%cit,
# J:Q = M_ETR, R:Y = M_PL). Inside each block, columns 3-4 of the block used
# to hold "IMF" data and columns 5-6 held "OECD (20%)" data:
#   col1,2 = GFA - Sales / GFA - Sales + Emp            (unchanged)
#   col3,4 = IMF - Sales / IMF - Sales + Emp             (old)
#   col5,6 = OECD (20%) - Sales / OECD (20%) - Sales + Emp (old)
#   col7,8 = OECD - Sales / OECD - Sales + Emp           (unchanged)
#
# A new "IMF (20%)" column pair (20% of the IMF figures, i.e. the
# stateless/offshore slice) is introduced right after GFA, the old IMF
# figures shift right into the slot that used to hold "OECD (20%)", and the
# old "OECD (20%)" figures are discarded. So after the edit:
#   col1,2 = GFA - Sales / GFA - Sales + Emp             (unchanged)
#   col3,4 = IMF (20%) - Sales / IMF (20%) - Sales + Emp (new = 0.2 * old IMF)
#   col5,6 = IMF - Sales / IMF - Sales + Emp             (= old col3,4)
#   col7,8 = OECD - Sales / OECD - Sales + Emp           (unchanged)
#
# Only the first block (B:I) actually has its IMF figures freshly derived
# (col3,4 <- 0.2 * old col3,4). The other two blocks (J:Q and R:Y) keep
# their original col3,4 values untouched and col5,6 is simply overwritten
# with a copy of col3,4 (matching the source workbook's own data, quirky
# as that is for those two blocks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each group: (GFA1, GFA2, IMF20pct-1, IMF20pct-2, IMF-1, IMF-2, OECD-1, OECD-2)
$groups = @(
    @{ Cols = @("B","C","D","E","F","G","H","I"); Scale = $true  },
    @{ Cols = @("J","K","L","M","N","O","P","Q"); Scale = $false },
    @{ Cols = @("R","S","T","U","V","W","X","Y"); Scale = $false }
)

# --- Row 2 header text: relabel the renumbered columns ---
foreach ($g in $groups) {
    $c = $g.Cols
    $ws.Range("$($c[2])2").Value = "IMF (20%) - Sales"
    $ws.Range("$($c[3])2").Value = "IMF (20%) - Sales + Emp"
    $ws.Range("$($c[4])2").Value = "IMF - Sales"
    $ws.Range("$($c[5])2").Value = "IMF - Sales + Emp"
}

# --- Data rows 4-10 ---
function Set-OrClear($range, $value) {
    if ($value -eq $null) {
        $range.ClearContents()
    } else {
        $range.Value = $value
    }
}

for ($row = 4; $row -le 10; $row++) {
    foreach ($g in $groups) {
        $c = $g.Cols

        $old3 = $ws.Range("$($c[2])$row").Value2
        $old4 = $ws.Range("$($c[3])$row").Value2

        if ($g.Scale) {
            if ($old3 -eq $null) { $new3 = $null } else { $new3 = $old3 * 0.2 }
            if ($old4 -eq $null) { $new4 = $null } else { $new4 = $old4 * 0.2 }
        } else {
            $new3 = $old3
            $new4 = $old4
        }
        $new5 = $old3
        $new6 = $old4

        Set-OrClear $ws.Range("$($c[4])$row") $new5
        Set-OrClear $ws.Range("$($c[5])$row") $new6
        Set-OrClear $ws.Range("$($c[2])$row") $new3
        Set-OrClear $ws.Range("$($c[3])$row") $new4
    }
}
